$d = $word.ActiveDocument

# 1) Simplify the headings: strip the "Objet TECHNICAL:TECHNICAL:" / "Type " prefixes
$d.Content.Find.Execute("Objet TECHNICAL:TECHNICAL:technical", $true, $false, $false, $false, $false, $true, 1, $false, "technical", 2)
$d.Content.Find.Execute("Type technicalObject", $true, $false, $false, $false, $false, $true, 1, $false, "technicalObject", 2)
$d.Content.Find.Execute("Type levelOneData", $true, $false, $false, $false, $false, $true, 1, $false, "levelOneData", 2)
$d.Content.Find.Execute("Type levelTwoData", $true, $false, $false, $false, $false, $true, 1, $false, "levelTwoData", 2)
$d.Content.Find.Execute("Type secondLevelTwoData", $true, $false, $false, $false, $false, $true, 1, $false, "secondLevelTwoData", 2)
$d.Content.Find.Execute("Type levelThreeData", $true, $false, $false, $false, $false, $true, 1, $false, "levelThreeData", 2)

# 2) Add a new row to the first table (end of table), describing a "nomenclatureField"
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "nomenclatureField"
$newRow.Cells.Item(2).Range.Text = "Nomenclature"
$newRow.Cells.Item(3).Range.Text = "string" + [char]11 + "(NOMENCLATURE: SI-SAMU-NOMENC_SEXE)"
$newRow.Cells.Item(4).Range.Text = "0..1"
$newRow.Cells.Item(5).Range.Text = "Enum from extenal nomenclature file"
# Leave the last ("Exemple") cell untouched so it stays an empty <w:r/> run
